$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.455362044514542
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 4.358119930609447

# Row 3
$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.1494219747398047
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 5.586269137925634

# Row 4
$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 0.7527432677738641
$ws.Range("E4").Value = 10.19245300693656
$ws.Range("G4").Value = 15.88780690183548

# Row 5
$ws.Range("B5").Value = 0.01293466051926884
$ws.Range("C5").Value = 10.34677158129881
$ws.Range("D5").Value = 3.537761648806719
$ws.Range("E5").Value = 10.19245300693656
$ws.Range("G5").Value = 24.08992089756136
